$d = $word.ActiveDocument

# --- Paragraph 1: merge the runs (no visible text change, but tidy) ---
# No action needed: text content is identical, only run-splitting changed,
# which is not observable via the Word OM / Find-Replace at the text level.

# --- Locate the first empty paragraph (between the first journal entry and
#     the trailing blank paragraph) and turn it into the "2024-05-02" bold
#     heading, then insert the new journal entry paragraph after it. ---

$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    if ($p.Range.Text.Trim().Length -eq 0) {
        $target = $p
        break
    }
}

$r = $target.Range
$r.Collapse(1)
$r.InsertAfter("2024-05-02")
$r.Font.Bold = $true

$newText = "Today, Sean was in the office, but he was on a call with a site that had a server failure, so he spent a lot of time working with them to get their server back in working order. He had me read some material about the historian software I'd be using to add reports from a different software so that the company that employed us could stop using it. In a nutshell the software is just making SQL queries from a database and displaying them nicely in a similar way to Excel. Afterwards we figured out how much access I had on the remote server housing the data I would need to access. Once we got access, Sean had me compile a list of all the important reports from the old software into an Excel spreadsheet so that they could be cross referenced once I started generating the new ones for the new software. Tomorrow, I'll be tuning in to the training session Sean is hosting for the employees that will be using the software so I can get some. experience using the software myself."

$endOfHeading = $target.Range
$endOfHeading.Collapse(0)
$endOfHeading.InsertParagraphAfter()
$newPara = $paras.Item($target.Index + 1)
$newRange = $newPara.Range
$newRange.Collapse(0)
$newRange.Font.Bold = $false
$newRange.InsertBefore($newText)
